$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Read original A:C values for rows 37..42
$vals = @()
for ($r = 37; $r -le 42; $r++) {
    $a = $ws.Cells.Item($r,1).Value2
    $b = $ws.Cells.Item($r,2).Value2
    $c = $ws.Cells.Item($r,3).Value2
    $vals += ,@($a,$b,$c)
}

# Write them shifted down by one row (37..42 -> 38..43)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $destRow = 38 + $i
    $ws.Cells.Item($destRow,1).Value = $vals[$i][0]
    $ws.Cells.Item($destRow,2).Value = $vals[$i][1]
    $ws.Cells.Item($destRow,3).Value = $vals[$i][2]
}

# Clear row 37 (A:C)
$ws.Cells.Item(37,1).Value = $null
$ws.Cells.Item(37,2).Value = $null
$ws.Cells.Item(37,3).Value = $null

for ($r = 37; $r -le 43; $r++) {
    $a = $ws.Cells.Item($r,1).Value2
    $b = $ws.Cells.Item($r,2).Value2
    $c = $ws.Cells.Item($r,3).Value2
    Write-Output ("Row " + $r + ": [" + $a + "] [" + $b + "] [" + $c + "]")
}
